$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data is stored as text (numbers-as-text), so format the
# target cells as Text before writing, keeping them consistent with the
# rest of the sheet instead of letting Excel auto-convert to numbers.
$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C4:E4").NumberFormat = "@"

# Swap the runs/balls/fours values between row 2 and row 4
$ws.Range("C2").Value = "11"
$ws.Range("D2").Value = "13"
$ws.Range("E2").Value = "1"

$ws.Range("C4").Value = "0"
$ws.Range("D4").Value = "1"
$ws.Range("E4").Value = "0"
